$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet: Correspond Handoff Datetime (E) and
# Correspond Handback DateTime (H) for rows 2 and 3.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 18:18:12"
$wsZh.Range("E3").Value = "2016-03-19 18:18:12"
$wsZh.Range("H2").Value = "2016-03-19 18:18:31"
$wsZh.Range("H3").Value = "2016-03-19 18:18:31"

# Update the "de-de" sheet: Correspond Handoff Datetime (E) and
# Correspond Handback DateTime (H) for rows 2 and 3.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 18:18:15"
$wsDe.Range("E3").Value = "2016-03-19 18:18:15"
$wsDe.Range("H2").Value = "2016-03-19 18:18:37"
$wsDe.Range("H3").Value = "2016-03-19 18:18:37"
